$wb = $excel.ActiveWorkbook

# Label the sheets with descriptive names (was Sheet1/Sheet2/Sheet3).
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Item(2)
$ws3 = $wb.Worksheets.Item(3)

$ws1.Name = "Table1-Prop_wt per Freq"
$ws2.Name = "Table2 Actual BMI vs perception"
$ws3.Name = "Table3-Proportion of Wt&Freq"

# Restore the selection on the first table sheet, then move the active
# tab/selection over to the second table sheet (Table2 ends up active).
$ws1.Activate()
$ws1.Range("B33").Select()

$ws2.Activate()
